# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new labels in AD1/AE1/AF1, formatted like the rest of
# the header row (bold, bordered, centered horizontally, top-aligned).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRng = $ws.Range("AD1:AF1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108   # xlCenter
$headerRng.VerticalAlignment = -4160     # xlTop
$headerRng.Borders.LineStyle = 1
$headerRng.Borders.Weight = 2

# Data rows 2-47: every row gets the same team record (79 wins, 83 losses,
# 0 ties).
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 79   # AD
    $ws.Cells.Item($r, 31).Value = 83   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
